$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting for price cells whose new values would otherwise be
# auto-converted to numbers by Excel (losing precision/trailing zeros/format).
$textCells = @("D5", "D6", "D11", "D14", "D20", "D21", "D22", "D24", "D26", "D29", "D31", "D32", "D38", "D41", "D42", "D45", "D47", "D49", "D50", "D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply updated cryptocurrency price and volume(1h) values
$ws.Range('D2').Value = '69.658.67'
$ws.Range('E2').Value = '  +4.89%  '
$ws.Range('D3').Value = '3.611.88'
$ws.Range('E3').Value = '  +4.84%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').Value = '631.08'
$ws.Range('E5').Value = '  +5.30%  '
$ws.Range('D6').Value = '158.70'
$ws.Range('E6').Value = '  +8.13%  '
$ws.Range('D7').Value = '3.609.17'
$ws.Range('E7').Value = '  +4.76%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('E9').Value = '  +4.28%  '
$ws.Range('E10').Value = '  +11.13%  '
$ws.Range('D11').Value = '7.57'
$ws.Range('E11').Value = '  +10.67%  '
$ws.Range('E12').Value = '  +6.78%  '
$ws.Range('E13').Value = '  +6.76%  '
$ws.Range('D14').Value = '33.76'
$ws.Range('E14').Value = '  +9.80%  '
$ws.Range('D15').Value = '4.226.77'
$ws.Range('E15').Value = '  +5.11%  '
$ws.Range('D16').Value = '3.618.78'
$ws.Range('E16').Value = '  +5.00%  '
$ws.Range('D17').Value = '69.605.60'
$ws.Range('E17').Value = '  +4.94%  '
$ws.Range('E18').Value = '  +1.31%  '
$ws.Range('E19').Value = '  +7.06%  '
$ws.Range('D20').Value = '16.17'
$ws.Range('E20').Value = '  +9.71%  '
$ws.Range('D21').Value = '10.26'
$ws.Range('E21').Value = '  +15.43%  '
$ws.Range('D22').Value = '463.45'
$ws.Range('E22').Value = '  +6.18%  '
$ws.Range('E23').Value = '  +5.15%  '
$ws.Range('D24').Value = '78.99'
$ws.Range('E24').Value = '  +3.28%  '
$ws.Range('E25').Value = '  +11.13%  '
$ws.Range('D26').Value = '10.76'
$ws.Range('E26').Value = '  +8.63%  '
$ws.Range('D27').Value = '3.758.81'
$ws.Range('E27').Value = '  +5.04%  '
$ws.Range('E28').Value = '  +0.01%  '
$ws.Range('D29').Value = '9.41'
$ws.Range('E29').Value = '  +15.61%  '
$ws.Range('E30').Value = '  +7.12%  '
$ws.Range('D31').Value = '1.73'
$ws.Range('E31').Value = '  +14.41%  '
$ws.Range('D32').Value = '0.172'
$ws.Range('E32').Value = '  +8.80%  '
$ws.Range('E33').Value = '  +9.21%  '
$ws.Range('E34').Value = '  -0.08%  '
$ws.Range('E35').Value = '  +7.30%  '
$ws.Range('E36').Value = '  +5.14%  '
$ws.Range('D37').Value = '3.613.07'
$ws.Range('E37').Value = '  +5.42%  '
$ws.Range('D38').Value = '8.48'
$ws.Range('E38').Value = '  +8.73%  '
$ws.Range('E39').Value = '  +15.55%  '
$ws.Range('E40').Value = '  +0.04%  '
$ws.Range('D41').Value = '0.0928'
$ws.Range('E41').Value = '  +8.98%  '
$ws.Range('D42').Value = '178.92'
$ws.Range('E42').Value = '  +3.47%  '
$ws.Range('E43').Value = '  +0.17%  '
$ws.Range('E44').Value = '  +6.73%  '
$ws.Range('D45').Value = '31.99'
$ws.Range('E45').Value = '  +24.64%  '
$ws.Range('E46').Value = '  +5.06%  '
$ws.Range('D47').Value = '1.39'
$ws.Range('E47').Value = '  +15.45%  '
$ws.Range('E48').Value = '  +13.56%  '
$ws.Range('D49').Value = '46.01'
$ws.Range('E49').Value = '  +1.73%  '
$ws.Range('D50').Value = '7.84'
$ws.Range('E50').Value = '  +4.98%  '
$ws.Range('D51').Value = '0.269'
